$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 208-209), shifting
# the existing rows 208-233 down to 210-235. This mirrors the diff, which
# prepends one more week's worth of "Frutilla" price data while leaving the
# rest of the historical rows untouched (just renumbered).
$ws.Rows("208:209").Insert()

# New row 208: Primera
$ws.Cells.Item(208, 1).Value = 7
$ws.Cells.Item(208, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(208, 3).Value = "Ñuble"
$ws.Cells.Item(208, 4).Value = 44578
$ws.Cells.Item(208, 5).Value = 16
$ws.Cells.Item(208, 6).Value = "Fruta"
$ws.Cells.Item(208, 7).Value = 100101
$ws.Cells.Item(208, 8).Value = "Berries"
$ws.Cells.Item(208, 9).Value = 100112025
$ws.Cells.Item(208, 10).Value = "Frutilla"
$ws.Cells.Item(208, 11).Value = "Sin especificar"
$ws.Cells.Item(208, 12).Value = "Primera"
$ws.Cells.Item(208, 13).Value = 300
$ws.Cells.Item(208, 14).Value = 6000
$ws.Cells.Item(208, 15).Value = 6500
$ws.Cells.Item(208, 16).Value = 6250
$ws.Cells.Item(208, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(208, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(208, 19).Value = 893
$ws.Cells.Item(208, 20).Value = 7

# New row 209: Segunda
$ws.Cells.Item(209, 1).Value = 7
$ws.Cells.Item(209, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(209, 3).Value = "Ñuble"
$ws.Cells.Item(209, 4).Value = 44578
$ws.Cells.Item(209, 5).Value = 16
$ws.Cells.Item(209, 6).Value = "Fruta"
$ws.Cells.Item(209, 7).Value = 100101
$ws.Cells.Item(209, 8).Value = "Berries"
$ws.Cells.Item(209, 9).Value = 100112025
$ws.Cells.Item(209, 10).Value = "Frutilla"
$ws.Cells.Item(209, 11).Value = "Sin especificar"
$ws.Cells.Item(209, 12).Value = "Segunda"
$ws.Cells.Item(209, 13).Value = 120
$ws.Cells.Item(209, 14).Value = 5000
$ws.Cells.Item(209, 15).Value = 5500
$ws.Cells.Item(209, 16).Value = 5250
$ws.Cells.Item(209, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(209, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(209, 19).Value = 750
$ws.Cells.Item(209, 20).Value = 7
